$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: Insert a new paragraph "Pomalejší než drát" (style "Odstavce 2.0")
# right before the "WAN" Heading 2 paragraph.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Nevýhodou je například", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $anchorPara = $rng.Paragraphs(1)
    $wanPara = $anchorPara.Next()
    $wanPara.Range.InsertParagraphBefore()
    $newPara = $anchorPara.Next()
    $newPara.Style = "Odstavce 2.0"
    $newPara.Range.Text = "Pomalejší než drát"
}

# ---------------------------------------------------------------------------
# Change 2: Mark the run that holds the first inline picture (the
# Obrázek 13 / anchorId 1E76D9E8 drawing, 4160520 x 3082925 EMU) as NoProof
# so it serializes with <w:rPr><w:noProof/></w:rPr>.
# ---------------------------------------------------------------------------
$shape = $d.InlineShapes.Item(1)
$shape.Range.NoProofing = 1

# ---------------------------------------------------------------------------
# Change 3: Merge the two runs "Kontrola," + " jestli soubor nebyl
# pozměněn " into a single run.
# ---------------------------------------------------------------------------
$rng = $d.Content
$null = $rng.Find.Execute("Kontrola, jestli soubor nebyl pozměněn ", $true, $false, $false, $false, $false, $true, 1, $false, "Kontrola, jestli soubor nebyl pozměněn ", 2)

# ---------------------------------------------------------------------------
# Change 4: Merge the runs " výkon" + " " (the 2nd & 3rd runs of the
# "Výpočetní výkon " paragraph) into a single run " výkon ", while leaving
# the separate "Výpočetní" run untouched.
#
# A plain Find/Replace across the whole paragraph coalesces *all* adjacent
# same-formatted runs (including the unrelated "Výpočetní" run), so instead
# we temporarily split the paragraph right after "Výpočetní", merge the
# runs in the now-isolated second paragraph, and rejoin the two paragraphs
# by deleting the temporary paragraph mark. That leaves "Výpočetní" as its
# own clean run and " výkon " as a single clean run, with no stray
# formatting left behind.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Výpočetní", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $splitPoint = $d.Range($rng.End, $rng.End)
    $splitPoint.InsertParagraphAfter()

    $para1 = $rng.Paragraphs(1)
    $para2 = $para1.Next()
    $null = $para2.Range.Find.Execute(" výkon ", $true, $false, $false, $false, $false, $true, 1, $false, " výkon ", 2)

    $para1b = $para2.Previous()
    $r96 = $para1b.Range
    $mark = $d.Range($r96.End - 1, $r96.End)
    $mark.Delete()
}

# ---------------------------------------------------------------------------
# Change 5: Merge the two runs "S" + "e odlišuje pro každou vrstvu v
# referenčním modelu ISO/OSI a obsahuje informace, které jsou relevantní
# pro tuto vrstvu. " into a single run.
# ---------------------------------------------------------------------------
$rng = $d.Content
$null = $rng.Find.Execute("Se odlišuje pro každou vrstvu v referenčním modelu ISO/OSI a obsahuje informace, které jsou relevantní pro tuto vrstvu. ", $true, $false, $false, $false, $false, $true, 1, $false, "Se odlišuje pro každou vrstvu v referenčním modelu ISO/OSI a obsahuje informace, které jsou relevantní pro tuto vrstvu. ", 2)
